$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.452439785003662
$ws.Range("B1").Value = 1.611706137657166
$ws.Range("C1").Value = 1.644370198249817
$ws.Range("D1").Value = 2.05788779258728
$ws.Range("E1").Value = 3.104574203491211
